$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the top of the data block (row 295),
# pushing all existing records from row 295..350 down by one row (to 296..351).
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the new weekly record.
$ws.Range("A295").Value = 10
$ws.Range("B295").Value = "Vega Modelo de Temuco"
$ws.Range("C295").Value = "La Araucanía"
$ws.Range("D295").Value = 44694
$ws.Range("E295").Value = 9
$ws.Range("F295").Value = 100112040
$ws.Range("G295").Value = "Cilantro"
$ws.Range("H295").Value = "Sin especificar"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 85
$ws.Range("K295").Value = 4000
$ws.Range("L295").Value = 5000
$ws.Range("M295").Value = 4588
$ws.Range("N295").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O295").Value = "Provincia de Cautín"
$ws.Range("P295").Value = 2294
$ws.Range("Q295").Value = 2
$ws.Range("R295").Value = "Hortaliza"
